$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Model description text (identical for rows 2-4)
$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       AdaBoostRegressor())]),`n                                            param_grid={'model__learning_rate': [0.1,`n                                                                                 0.5,`n                                                                                 1.0],`n                                                        'model__n_estimators': [50,`n                                                                                100,`n                                                                                150]},`n                                            scoring='neg_mean_squared_error'))"

# New header cell F1 - copy header style from E1, then set its value
$ws.Range("F1").Value = "Modelo"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Update numeric metrics for rows 2-4 (columns B, C, D)
$ws.Range("B2").Value = 0.0510073685469512
$ws.Range("C2").Value = 0.9984999603120929
$ws.Range("D2").Value = 0.1657582300642175

$ws.Range("B3").Value = 0.07962121786079235
$ws.Range("C3").Value = 0.9992447570665006
$ws.Range("D3").Value = 0.209259200321973

$ws.Range("B4").Value = 0.08780122833485468
$ws.Range("C4").Value = 0.998811225566325
$ws.Range("D4").Value = 0.2399893346931228

# New column F (Modelo) values for data rows
$ws.Range("F2").Value = $modelText
$ws.Range("F3").Value = $modelText
$ws.Range("F4").Value = $modelText

# The multi-line text triggers an automatic row-height expansion; restore
# the rows to their default auto-fit height so no spurious height is stored.
$ws.Rows.Item(2).AutoFit() | Out-Null
$ws.Rows.Item(3).AutoFit() | Out-Null
$ws.Rows.Item(4).AutoFit() | Out-Null

$excel.CutCopyMode = 0
